$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Plan1")

# --- Cell text content changes ---
# (order matches the order new shared strings were appended in the target file)

# B47-B49: Solicitante -> RepresentanteLegalUm
$ws.Range("B47").Value = "Nome/Procurador Legal:  `$RepresentanteLegalUm"
$ws.Range("B48").Value = "Telefone: `$TelefoneCelularRepresentanteLegalUm"
$ws.Range("B49").Value = "E-mail: `$EmailRepresentanteLegalUm"

# G7: Nº:  $NumeroEnel -> Nº:  $NrEnel
$ws.Range("G7").Value = "Nº:  `$NrEnel"

# B13: merge in the coordinates placeholder, and clear G13
$ws.Range("B13").Value = "Localização em coordenadas ( Latitude , Longitude) :  `$CoordenadasGD "
$ws.Range("G13").ClearContents()

# --- Column width change: column F narrower (target OOXML width 7.140625 chars) ---
$ws.Columns.Item(6).ColumnWidth = 6.25

# --- Sheet view changes: selection ---
$ws.Activate()
$ws.Range("I11").Select()

$wb.Save()
